$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.661.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.17%  "
$ws.Range("D3").Value = "'1.891.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'244.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.4966"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.2959"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").Value = "'0.06811"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.44%  "
$ws.Range("D10").Value = "'1.890.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").Value = "'17.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("D12").Value = "'0.07310"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").Value = "'90.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.50%  "
$ws.Range("D14").Value = "'5.070"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.85%  "
$ws.Range("D15").Value = "'0.6722"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "'30.640.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").Value = "'0.000007928"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "'13.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.98%  "
$ws.Range("D20").Value = "'2.136.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "'4.851"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("D23").Value = "'176.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +31.20%  "
$ws.Range("D24").Value = "'6.060"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.44%  "
$ws.Range("D25").Value = "'9.274"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("D26").Value = "'154.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.10%  "
$ws.Range("D27").Value = "'18.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.25%  "
$ws.Range("D28").Value = "'1.925"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").Value = "'1.391"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("D30").Value = "'4.329"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.10%  "
$ws.Range("D31").Value = "'0.08921"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.81%  "
$ws.Range("D32").Value = "'4.027"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("D33").Value = "'0.05223"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.35%  "
$ws.Range("D34").Value = "'0.7394"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.20%  "
$ws.Range("D35").Value = "'1.135"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.50%  "
$ws.Range("D36").Value = "'2.685"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").Value = "'0.01870"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.37%  "
$ws.Range("D38").Value = "'2.703"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'2.170"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'0.9355"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").Value = "'0.4362"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.34%  "
$ws.Range("D42").Value = "'105.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.89%  "
$ws.Range("D43").Value = "'5.812"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").Value = "'7.658"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.48%  "
$ws.Range("E46").Value = "  +7.96%  "
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("D48").Value = "'33.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.23%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.521"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.17%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").Value = "'0.3876"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.09%  "
$ws.Range("D51").Value = "'1.379"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.46%  "
